$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for team record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold/centered/bordered header style used by the rest of row 1
$hdr = $ws.Range("AD1:AF1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160

# Fill team record values (Wins=70, Losses=92, Ties=0) for every player row
$lastRow = 56

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 70
    $ws.Cells.Item($r, 31).Value = 92
    $ws.Cells.Item($r, 32).Value = 0
}

$ws.Range("A1").Select()
